$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2024-04-13)
$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 6.201049113329182

# Row 3 (2024-03-30)
$ws.Range("B3").Value = 1.459612070389937
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 29.84159230404497
